# missions.xlsx update
# - re-stamp the DateDebut/DateFin (D/E) columns of the existing mission rows onto a new
#   "yyyy-mm-dd" date style
# - normalise row 13 (previously a plain-text date) into a real date value using that style
# - append the newly logged mission as row 14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# ---------------------------------------------------------------------------
# 1) DateDebut / DateFin for every existing mission row (2-12): move onto a
#    fresh date style.
# ---------------------------------------------------------------------------
$ws.Range("D2:E12").NumberFormat = "yyyy-mm-dd"

# ---------------------------------------------------------------------------
# 2) Row 13 still held its DateDebut/DateFin as plain text ("2024-08-22").
#    Replace with real date serials using the same style as the rest of the
#    table.
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 4).Value = 45526
$ws.Cells.Item(13, 5).Value = 45526
$ws.Range("D13:E13").NumberFormat = "yyyy-mm-dd"

# ---------------------------------------------------------------------------
# 3) Append the newly logged mission (row 14).
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = 0
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = "20:20:04"

# DateDebut/DateFin on a freshly appended row are still plain text at this
# point (only later normalised, as row 13 was above). Build the text through
# a TEXT() formula and flatten it to a static value so Excel's "this looks
# like a date" auto-detection doesn't silently reformat the cell the way a
# plain string assignment would.
$ws.Cells.Item(14, 4).Formula = '=TEXT("2024-08-26","yyyy-mm-dd")'
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial($xlPasteValues)

$ws.Cells.Item(14, 5).Formula = '=TEXT("2024-08-26","yyyy-mm-dd")'
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(14, 5).PasteSpecial($xlPasteValues)

$ws.Cells.Item(14, 6).Value = "Cartographie"
$ws.Cells.Item(14, 7).Value = "Ete"
$ws.Cells.Item(14, 8).Value = "DroneA"
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 3
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 1
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0

$excel.CutCopyMode = $false
